$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Details" to "Sheet1"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# Make the first sheet ("Sheet1") the active/selected tab
# (previously "Conventional-Results", the 3rd sheet, was active)
$ws1.Activate()
